# Update values on the "constants" sheet to reflect new uncertainty-run
# progress graph results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

$ws.Range("B2").Value = 9.817290327844651
$ws.Range("B3").Value = 0.16
$ws.Range("B4").Value = 1849.448887606751
$ws.Range("B5").Value = 24267.59167219321
$ws.Range("B8").Value = 0.7705065931625953
$ws.Range("B9").Value = 0.6994721573256851
$ws.Range("B10").Value = 2.894718404854602
$ws.Range("B11").Value = 0.3826378407848143
